$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (B2, C2, C3, B4, C5)
$ws.Range("B2").Value = 5.3
$ws.Range("C2").Value = 10.8
$ws.Range("C3").Value = 9.7
$ws.Range("B4").Value = 0.75
$ws.Range("C5").Value = 21

# Update the active selection to C8
$ws.Range("C8").Select()

# Update window size
$excel.ActiveWindow.Width = 22785
$excel.ActiveWindow.Height = 10575
